# Applies updated "dSF" (column F) values for the 2022 bundy_dylan save_data sheet.
# These values represent re-pulled / recomputed delta-score-final figures
# (per commit message: "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    5  = 1
    6  = -7
    7  = 2
    8  = -1
    9  = -1
    10 = 2
    11 = 2
    13 = 1
    14 = -2
    15 = 5
    16 = 0
    17 = 1
    18 = -3
    19 = -1
    20 = 10
    21 = 4
    22 = -4
    24 = 1
    26 = -4
    27 = -5
    28 = 7
    29 = 3
    30 = 1
    31 = 1
    32 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
